$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3544.9
$ws.Range("I40").Value = 2778.4285
$ws.Range("K40").Value = 2778.4285
$ws.Range("M40").Value = -2603.4285
$ws.Range("H43").Value = 10000.333
$ws.Range("I43").Value = 10000.5
$ws.Range("K43").Value = 10000.5
$ws.Range("M43").Value = -9931.5
$ws.Range("H62").Value = 6233.6206
$ws.Range("I62").Value = 5621.174
$ws.Range("J62").Value = 8581.333000000001
$ws.Range("K62").Value = 5621.174
$ws.Range("L62").Value = 8581.333000000001
$ws.Range("M62").Value = -4997.174
$ws.Range("N62").Value = -9829.333000000001
$ws.Range("H65").Value = 6233.6206
$ws.Range("I65").Value = 5621.174
$ws.Range("J65").Value = 8581.333000000001
$ws.Range("K65").Value = 28105.87
$ws.Range("L65").Value = 42906.665
$ws.Range("M65").Value = -24985.87
$ws.Range("N65").Value = -49146.665
$ws.Range("H132").Value = 1968239.5
$ws.Range("I132").Value = 2210107.2
$ws.Range("K132").Value = 6630321.600000001
$ws.Range("M132").Value = -6627791.600000001
$ws.Range("H137").Value = 25080.268
$ws.Range("I137").Value = 37656.223
$ws.Range("K137").Value = 112968.669
$ws.Range("M137").Value = -110418.669
$ws.Range("H141").Value = 1497.8823
$ws.Range("I141").Value = 1132.6666
$ws.Range("K141").Value = 3397.9998
$ws.Range("M141").Value = 1782.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16429.732
$ws.Range("I32").Value = 17260.457
$ws.Range("J32").Value = 4799.6
$ws.Range("K32").Value = 17260.457
$ws.Range("L32").Value = 4799.6
$ws.Range("M32").Value = -16973.457
$ws.Range("N32").Value = -5373.6
$ws.Range("H45").Value = 3839.5789
$ws.Range("I45").Value = 1735.75
$ws.Range("J45").Value = 7446.143
$ws.Range("K45").Value = 1735.75
$ws.Range("L45").Value = 7446.143
$ws.Range("M45").Value = -1358.75
$ws.Range("N45").Value = -8200.143
$ws.Range("H61").Value = 4469.6562
$ws.Range("I61").Value = 1007.6316
$ws.Range("K61").Value = 1007.6316
$ws.Range("M61").Value = -795.6316
$ws.Range("H88").Value = 10983.167
$ws.Range("I88").Value = 2999
$ws.Range("J88").Value = 12580
$ws.Range("K88").Value = 2999
$ws.Range("L88").Value = 12580
$ws.Range("M88").Value = -2593
$ws.Range("N88").Value = -13392
$ws.Range("H91").Value = 10983.167
$ws.Range("I91").Value = 2999
$ws.Range("J91").Value = 12580
$ws.Range("K91").Value = 2999
$ws.Range("L91").Value = 12580
$ws.Range("M91").Value = -1595
$ws.Range("N91").Value = -15388
$ws.Range("H132").Value = 1351.2319
$ws.Range("I132").Value = 1087.6383
$ws.Range("J132").Value = 1914.3636
$ws.Range("K132").Value = 3262.9149
$ws.Range("L132").Value = 5743.0908
$ws.Range("M132").Value = -732.9149000000002
$ws.Range("N132").Value = -10803.0908
$ws.Range("H136").Value = 4469.6562
$ws.Range("I136").Value = 1007.6316
$ws.Range("K136").Value = 3022.8948
$ws.Range("M136").Value = -472.8948

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1765.9131
$ws.Range("I86").Value = 1740.0588
$ws.Range("J86").Value = 1839.1666
$ws.Range("K86").Value = 1740.0588
$ws.Range("L86").Value = 1839.1666
$ws.Range("M86").Value = -617.0588
$ws.Range("N86").Value = -4085.1666
$ws.Range("H89").Value = 1765.9131
$ws.Range("I89").Value = 1740.0588
$ws.Range("J89").Value = 1839.1666
$ws.Range("K89").Value = 8700.294
$ws.Range("L89").Value = 9195.833000000001
$ws.Range("M89").Value = -3084.294
$ws.Range("N89").Value = -20427.833
$ws.Range("H134").Value = 3607.0732
$ws.Range("I134").Value = 1862.6154
$ws.Range("K134").Value = 5587.8462
$ws.Range("M134").Value = -3052.8462

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 25774.426
$ws.Range("I132").Value = 30037.617
$ws.Range("K132").Value = 90112.851
$ws.Range("M132").Value = -87582.851
$ws.Range("H135").Value = 119996.75
$ws.Range("J135").Value = 119996.75
$ws.Range("L135").Value = 119996.75
$ws.Range("N135").Value = -130136.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 2626.1
$ws.Range("I137").Value = 2362.3333
$ws.Range("K137").Value = 7086.999899999999
$ws.Range("M137").Value = -1986.999899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 755.3333
$ws.Range("I2").Value = 864
$ws.Range("J2").Value = 584.5714
$ws.Range("K2").Value = 864
$ws.Range("L2").Value = 584.5714
$ws.Range("M2").Value = -751
$ws.Range("N2").Value = -810.5714
$ws.Range("H102").Value = 18290.695
$ws.Range("I102").Value = 21527.967
$ws.Range("J102").Value = 2104.3333
$ws.Range("K102").Value = 21527.967
$ws.Range("L102").Value = 2104.3333
$ws.Range("M102").Value = -19905.967
$ws.Range("N102").Value = -5348.3333
$ws.Range("H126").Value = 3456.25
$ws.Range("I126").Value = 1599.625
$ws.Range("J126").Value = 4694
$ws.Range("K126").Value = 4798.875
$ws.Range("L126").Value = 14082
$ws.Range("M126").Value = -2328.875
$ws.Range("N126").Value = -19022
$ws.Range("H132").Value = 2227.7856
$ws.Range("I132").Value = 2233.8076
$ws.Range("J132").Value = 2149.5
$ws.Range("K132").Value = 6701.4228
$ws.Range("L132").Value = 6448.5
$ws.Range("M132").Value = -4171.4228
$ws.Range("N132").Value = -11508.5
$ws.Range("H139").Value = 112712
$ws.Range("J139").Value = 112712
$ws.Range("L139").Value = 112712
$ws.Range("N139").Value = -122992

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2206.5518
$ws.Range("I7").Value = 1835.9546
$ws.Range("K7").Value = 1835.9546
$ws.Range("M7").Value = -1723.9546
$ws.Range("H22").Value = 1272.2222
$ws.Range("I22").Value = 740
$ws.Range("J22").Value = 1937.5
$ws.Range("K22").Value = 740
$ws.Range("L22").Value = 1937.5
$ws.Range("M22").Value = -445
$ws.Range("N22").Value = -2527.5
$ws.Range("H27").Value = 1272.2222
$ws.Range("I27").Value = 740
$ws.Range("J27").Value = 1937.5
$ws.Range("K27").Value = 740
$ws.Range("L27").Value = 1937.5
$ws.Range("M27").Value = -633
$ws.Range("N27").Value = -2151.5
$ws.Range("H40").Value = 2490.55
$ws.Range("I40").Value = 2400.353
$ws.Range("K40").Value = 2400.353
$ws.Range("M40").Value = -2264.353
$ws.Range("H61").Value = 757.5172
$ws.Range("I61").Value = 607.375
$ws.Range("K61").Value = 607.375
$ws.Range("M61").Value = -405.375
$ws.Range("H113").Value = 757.5172
$ws.Range("I113").Value = 607.375
$ws.Range("K113").Value = 607.375
$ws.Range("M113").Value = 1562.625
$ws.Range("H126").Value = 2206.5518
$ws.Range("I126").Value = 1835.9546
$ws.Range("K126").Value = 5507.8638
$ws.Range("M126").Value = -3037.8638
$ws.Range("H132").Value = 5163.1333
$ws.Range("I132").Value = 5035.769
$ws.Range("J132").Value = 5991
$ws.Range("K132").Value = 15107.307
$ws.Range("L132").Value = 17973
$ws.Range("M132").Value = -12577.307
$ws.Range("N132").Value = -23033
$ws.Range("H134").Value = 68497.5
$ws.Range("J134").Value = 68497.5
$ws.Range("L134").Value = 68497.5
$ws.Range("N134").Value = -78637.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 556966.6
$ws.Range("I4").Value = 819.9091
$ws.Range("K4").Value = 819.9091
$ws.Range("M4").Value = -706.9091
$ws.Range("H126").Value = 168646.6
$ws.Range("I126").Value = 1661.9524
$ws.Range("K126").Value = 4985.857199999999
$ws.Range("M126").Value = -2515.857199999999
$ws.Range("H136").Value = 13180.327
$ws.Range("I136").Value = 15098.2
$ws.Range("J136").Value = 4549.9
$ws.Range("K136").Value = 45294.60000000001
$ws.Range("L136").Value = 13649.7
$ws.Range("M136").Value = -42744.60000000001
$ws.Range("N136").Value = -18749.7
$ws.Range("H137").Value = 97165.55499999999
$ws.Range("J137").Value = 108248.336
$ws.Range("L137").Value = 108248.336
$ws.Range("N137").Value = -118448.336
$ws.Range("H141").Value = 97552.86
$ws.Range("J141").Value = 97552.86
$ws.Range("L141").Value = 97552.86
$ws.Range("N141").Value = -107912.86
